# Auto-generated Excel COM-interop script to apply numeric updates
# to the Leve profit-tracking tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row the scheduled price-refresh run recalculated the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N); some rows
# gained or lost a LeveProfitNQ/LeveProfitHQ cell depending on whether that
# side of the recipe now has data.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 117.3125
$ws.Range("I6").Value = 123.46667
$ws.Range("K6").Value = 370.40001
$ws.Range("M6").Value = -258.40001

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H31").Value = 189.75
$ws.Range("I31").Value = 51.666668
$ws.Range("J31").Value = 604
$ws.Range("K31").Value = 155.000004
$ws.Range("L31").Value = 1812
$ws.Range("M31").Value = 74.99999600000001
$ws.Range("N31").Value = -2272

$ws.Range("H99").Value = 1603
$ws.Range("I99").Value = 1603
$ws.Range("K99").Value = 4809
$ws.Range("M99").Value = -3311

$ws.Range("H135").Value = 768.8125
$ws.Range("I135").Value = 901.2308
$ws.Range("J135").Value = 195
$ws.Range("K135").Value = 8111.077200000001
$ws.Range("L135").Value = 1755
$ws.Range("M135").Value = -5576.077200000001
$ws.Range("N135").Value = -6825

$ws.Range("H137").Value = 2223.7856
$ws.Range("I137").Value = 1682.5555
$ws.Range("J137").Value = 3198
$ws.Range("K137").Value = 5047.666499999999
$ws.Range("L137").Value = 9594
$ws.Range("M137").Value = -2497.666499999999
$ws.Range("N137").Value = -14694

$ws.Range("H138").Value = 13222.5
$ws.Range("J138").Value = 13256.516
$ws.Range("L138").Value = 39769.548
$ws.Range("N138").Value = -50049.548

$ws.Range("H141").Value = 1231.3334
$ws.Range("I141").Value = 1231.3334
$ws.Range("K141").Value = 3694.0002
$ws.Range("M141").Value = 1485.9998


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10826.849
$ws.Range("I32").Value = 9912.483
$ws.Range("K32").Value = 9912.483
$ws.Range("M32").Value = -9625.483

$ws.Range("H61").Value = 4656.857
$ws.Range("I61").Value = 2533.3333
$ws.Range("K61").Value = 2533.3333
$ws.Range("M61").Value = -2321.3333

$ws.Range("H88").Value = 3200.5
$ws.Range("I88").Value = 2266.5
$ws.Range("J88").Value = 4134.5
$ws.Range("K88").Value = 2266.5
$ws.Range("L88").Value = 4134.5
$ws.Range("M88").Value = -1860.5
$ws.Range("N88").Value = -4946.5

$ws.Range("H91").Value = 3200.5
$ws.Range("I91").Value = 2266.5
$ws.Range("J91").Value = 4134.5
$ws.Range("K91").Value = 2266.5
$ws.Range("L91").Value = 4134.5
$ws.Range("M91").Value = -862.5
$ws.Range("N91").Value = -6942.5

$ws.Range("H97").Value = 700
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 700
$ws.Range("M97").Value = -204

$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122

$ws.Range("H132").Value = 2015.5264
$ws.Range("I132").Value = 1618.8462
$ws.Range("K132").Value = 4856.5386
$ws.Range("M132").Value = -2326.5386

$ws.Range("H136").Value = 4656.857
$ws.Range("I136").Value = 2533.3333
$ws.Range("K136").Value = 7599.999899999999
$ws.Range("M136").Value = -5049.999899999999


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2776.1304
$ws.Range("I20").Value = 2063.2856
$ws.Range("J20").Value = 3885
$ws.Range("K20").Value = 2063.2856
$ws.Range("L20").Value = 3885
$ws.Range("M20").Value = -1816.2856
$ws.Range("N20").Value = -4379

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H94").Value = 2280.125
$ws.Range("I94").Value = 1605.9286
$ws.Range("K94").Value = 1605.9286
$ws.Range("M94").Value = -1154.9286

$ws.Range("H99").Value = 1548.3572
$ws.Range("I99").Value = 1379.2
$ws.Range("J99").Value = 1971.25
$ws.Range("K99").Value = 1379.2
$ws.Range("L99").Value = 1971.25
$ws.Range("M99").Value = 118.8
$ws.Range("N99").Value = -4967.25

$ws.Range("H102").Value = 35000
$ws.Range("I102").Value = 30000
$ws.Range("J102").Value = 40000
$ws.Range("K102").Value = 30000
$ws.Range("L102").Value = 40000
$ws.Range("M102").Value = -26755
$ws.Range("N102").Value = -46490

$ws.Range("H105").Value = 5090.909
$ws.Range("I105").Value = 4285.7144
$ws.Range("K105").Value = 4285.7144
$ws.Range("M105").Value = -2538.7144

$ws.Range("H134").Value = 2650
$ws.Range("I134").Value = 2200
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6600
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4065
$ws.Range("N134").Value = -17070


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4990.8335
$ws.Range("I86").Value = 3992.6
$ws.Range("K86").Value = 3992.6
$ws.Range("M86").Value = -2869.6

$ws.Range("H89").Value = 4990.8335
$ws.Range("I89").Value = 3992.6
$ws.Range("K89").Value = 19963
$ws.Range("M89").Value = -14347

$ws.Range("H134").Value = 3744.3704
$ws.Range("J134").Value = 5011.7144
$ws.Range("L134").Value = 15035.1432
$ws.Range("N134").Value = -20105.1432


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 3131.9333
$ws.Range("I86").Value = 999.3333
$ws.Range("J86").Value = 3665.0833
$ws.Range("K86").Value = 2997.9999
$ws.Range("L86").Value = 10995.2499
$ws.Range("M86").Value = -1811.9999
$ws.Range("N86").Value = -13367.2499

$ws.Range("H88").Value = 8000
$ws.Range("J88").Value = 8000
$ws.Range("L88").Value = 24000
$ws.Range("N88").Value = -24856

$ws.Range("H89").Value = 3131.9333
$ws.Range("I89").Value = 999.3333
$ws.Range("J89").Value = 3665.0833
$ws.Range("K89").Value = 8993.9997
$ws.Range("L89").Value = 32985.7497
$ws.Range("M89").Value = -3065.9997
$ws.Range("N89").Value = -44841.7497

$ws.Range("H91").Value = 8000
$ws.Range("J91").Value = 8000
$ws.Range("L91").Value = 24000
$ws.Range("N91").Value = -26964

$ws.Range("H109").Value = 1339
$ws.Range("I109").Value = 897.5
$ws.Range("J109").Value = 1633.3334
$ws.Range("K109").Value = 2692.5
$ws.Range("L109").Value = 4900.0002
$ws.Range("M109").Value = -1652.5
$ws.Range("N109").Value = -6980.0002

$ws.Range("H131").Value = 1739.4
$ws.Range("I131").Value = 1599
$ws.Range("J131").Value = 1774.5
$ws.Range("K131").Value = 4797
$ws.Range("L131").Value = 5323.5
$ws.Range("M131").Value = 243
$ws.Range("N131").Value = -15403.5

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25346

$ws.Range("H80").Value = 2663.7273
$ws.Range("I80").Value = 1514.4286
$ws.Range("J80").Value = 4675
$ws.Range("K80").Value = 1514.4286
$ws.Range("L80").Value = 4675
$ws.Range("M80").Value = -516.4286
$ws.Range("N80").Value = -6671

$ws.Range("H83").Value = 2663.7273
$ws.Range("I83").Value = 1514.4286
$ws.Range("J83").Value = 4675
$ws.Range("K83").Value = 7572.143
$ws.Range("L83").Value = 23375
$ws.Range("M83").Value = -2580.143
$ws.Range("N83").Value = -33359

$ws.Range("H113").Value = 1074.1111
$ws.Range("I113").Value = 1074.1111
$ws.Range("K113").Value = 1074.1111
$ws.Range("M113").Value = 1095.8889


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1291.8334
$ws.Range("I22").Value = 988.8889
$ws.Range("J22").Value = 2200.6667
$ws.Range("K22").Value = 988.8889
$ws.Range("L22").Value = 2200.6667
$ws.Range("M22").Value = -693.8889
$ws.Range("N22").Value = -2790.6667

$ws.Range("H27").Value = 1291.8334
$ws.Range("I27").Value = 988.8889
$ws.Range("J27").Value = 2200.6667
$ws.Range("K27").Value = 988.8889
$ws.Range("L27").Value = 2200.6667
$ws.Range("M27").Value = -881.8889
$ws.Range("N27").Value = -2414.6667

$ws.Range("H55").Value = 953.0625
$ws.Range("I55").Value = 356.7143
$ws.Range("J55").Value = 1416.8889
$ws.Range("K55").Value = 356.7143
$ws.Range("L55").Value = 1416.8889
$ws.Range("M55").Value = -183.7143
$ws.Range("N55").Value = -1762.8889

$ws.Range("H132").Value = 4001.625
$ws.Range("I132").Value = 2002.1666
$ws.Range("K132").Value = 6006.4998
$ws.Range("M132").Value = -3476.4998


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H75").Value = 24974.75
$ws.Range("I75").Value = 24949.5
$ws.Range("K75").Value = 24949.5
$ws.Range("M75").Value = -24013.5

$ws.Range("H78").Value = 24974.75
$ws.Range("I78").Value = 24949.5
$ws.Range("K78").Value = 74848.5
$ws.Range("M78").Value = -70168.5

$ws.Range("H96").Value = 2039.9
$ws.Range("I96").Value = 1612.5
$ws.Range("K96").Value = 1612.5
$ws.Range("M96").Value = -239.5

$ws.Range("H126").Value = 1954.6
$ws.Range("I126").Value = 1443.375
$ws.Range("K126").Value = 4330.125
$ws.Range("M126").Value = -1860.125

